$d = $word.ActiveDocument

# Helper: force a run-split at a collapsed point by adding and
# immediately deleting a throwaway bookmark there. This mirrors how
# Word splits runs at the caret when you type/insert/format text,
# without leaving any residual formatting behind.
function Split-At($pos) {
    $pt = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TempSplit", $pt)
    $d.Bookmarks("TempSplit").Delete()
}

# ============================================================
# Change 2 (Android Development bullet): heal the stray mid-word
# run split ("Android-spec" | "ific design issues.") back into a
# single run, and drop the _GoBack bookmark that sat there (it
# relocates to the edit made below).
# ============================================================
$d.Content.Find.Execute( `
    "Cards in this area focus on building Android layout elements, building logic to populate the layout with data, and working on Android-specific design issues.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Cards in this area focus on building Android layout elements, building logic to populate the layout with data, and working on Android-specific design issues.", `
    2)

# ============================================================
# Change 1: Customer Description paragraph updates
#   "...service ticket application should be created..." ->
#   "...service ticket application will be created..."
#   "...service technicians. Innovative needs..." ->
#   "...service technicians. Innovative Systems needs..."
# ============================================================

# Locate "should" (within "should be created") and capture its start
# position before editing it.
$rngShould = $d.Content
$rngShould.Find.Execute("should be created", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$shouldStart = $rngShould.Start

# Replace just the word "should" (6 chars) with "will" (4 chars).
$wordOnly = $d.Range($shouldStart, $shouldStart + 6)
$wordOnly.Text = "will"
$willStart = $shouldStart
$willEnd = $shouldStart + 4

# Locate "Innovative needs" and insert " Systems" right after
# "Innovative" (10 chars in), before the existing " needs".
$rng2 = $d.Content
$rng2.Find.Execute("Innovative needs", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$innovEnd = $rng2.Start + 10
$insertPt = $d.Range($innovEnd, $innovEnd)
$insertPt.InsertAfter(" Systems")
$sysStart = $innovEnd
$sysMid = $innovEnd + 8   # right after the inserted " Systems"

# Split runs at exactly the boundaries Word itself would have left
# behind from these two small edits.
Split-At $willStart
Split-At $willEnd
Split-At $sysStart
Split-At $sysMid

# Move the _GoBack bookmark (Word's "last edit" marker) to this new
# edit location, right after the inserted " Systems".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmPt = $d.Range($sysMid, $sysMid)
$d.Bookmarks.Add("_GoBack", $bmPt)
